$d = $word.ActiveDocument

# Locate the "!!" at the end of "Let's get started!!" via Find so we don't
# depend on hard-coded character offsets.
$bangRng = $d.Content.Duplicate
$found = $bangRng.Find.Execute("!!", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the target '!!' run to edit."
}

$para = $bangRng.Paragraphs(1)

# Remove the existing "!!" run outright -- we are going to re-add its text
# (plus the new surrounding bold runs) as fresh runs appended at the very
# end of the paragraph, which is the only insertion point that reliably
# yields discrete <w:r> elements instead of being coalesced into whatever
# run happens to sit to the left of the insertion point.
$bangRng.Delete()

function Append-BoldRun($doc, $paragraph, $text) {
    $endPos = $paragraph.Range.End - 1
    $insertion = $doc.Range($endPos, $endPos)
    $insertion.InsertAfter($text)
    $newRun = $doc.Range($endPos, $endPos + $text.Length)
    # Toggling Bold off/on (even though the end value is unchanged) forces
    # the engine to treat this span as its own run rather than merging it
    # into a neighboring run that already happens to be bold.
    $newRun.Font.Bold = $false
    $newRun.Font.Bold = $true
}

Append-BoldRun $d $para " "
Append-BoldRun $d $para " "
Append-BoldRun $d $para "!!"
Append-BoldRun $d $para "!"

Write-Output $para.Range.Text
